# Tidy up BDM info: remove the Evan / Simon / Kriti / Jack / pseudo-Sally
# questions, keep the Mishika / Christian / Sabi / Borys / lottery-Sally
# questions, and add a "Both of the above options" answer choice to the
# Sally lottery question.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows from the bottom up so row numbers of not-yet-processed rows
# stay stable.
$ws.Rows("11").Delete()   # Jack spoon/fork question
$ws.Rows("6").Delete()    # Kriti mug question
$ws.Rows("5").Delete()    # Simon ice-cream question
$ws.Rows("3").Delete()    # Sally "pseudo switch-point" question (duplicate)
$ws.Rows("2").Delete()    # Evan cheeseburger question

# The remaining Sally question (previously row 4) is now row 2. Re-jig its
# answer options: add a new "Both of the above options" choice and reorder
# the existing ones to match.
$ws.Range("B2").Value2 = "Sally will spend more time in the experiment"
$ws.Range("C2").Value2 = "Sally may not always get her preferred outcome"
$ws.Range("D2").Value2 = "Both of the above options"
$ws.Range("E2").Value2 = "There are no negative consequences associated with lying for Sally"
$ws.Range("F2").Value2 = "B"

$ws.Range("C2").Select()
